$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# --- Row 87: new date entry ---
# Copy date formatting (style) from an existing date cell (A84) so we reuse
# the existing style index instead of creating a new number format.
$ws.Cells.Item(84,1).Copy()
$ws.Cells.Item(87,1).PasteSpecial(-4122)
$ws.Cells.Item(87,1).Value2 = 44585

# Copy hour-count cell formatting (style) from an existing B cell.
$ws.Cells.Item(86,2).Copy()
$ws.Cells.Item(87,2).PasteSpecial(-4122)
$ws.Cells.Item(87,2).Value2 = 1

$ws.Cells.Item(87,3).Value = "analyzer form info, home info kirjoitusta"
$ws.Cells.Item(87,4).Value = "client"

# --- Row 88 ---
$ws.Cells.Item(86,2).Copy()
$ws.Cells.Item(88,2).PasteSpecial(-4122)
$ws.Cells.Item(88,2).Value2 = 3

$ws.Cells.Item(88,3).Value = "DataView Refaktorointi, 3 uutta reduceria, componentteja siistitty niin että logiikkaa siirretty pois komponentista"
$ws.Cells.Item(88,4).Value = "client"

# --- Row 89 ---
$ws.Cells.Item(86,2).Copy()
$ws.Cells.Item(89,2).PasteSpecial(-4122)
$ws.Cells.Item(89,2).Value2 = 1

$ws.Cells.Item(89,3).Value = "Market css refaktorointia, toimintojen testausta, kaikki toimii muutosten jälkeen"
$ws.Cells.Item(89,4).Value = "client"

# --- Update summary formulas to include the new rows ---
$ws.Range("B96").Formula = "=SUM(B2:B89)"
$ws.Range("B98").Formula = "=B96/B97*100"

$excel.CalculateFull()

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("D89").Select()

$wb.Save()
